$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.9444154798013937
$ws.Range("E2").Value = 0.9444154798013937

# Row 3
$ws.Range("D3").Value = 0.0003106111345444427
$ws.Range("E3").Value = 0.0003106111345444427

# Row 4
$ws.Range("D4").Value = 0.003128754952689148
$ws.Range("E4").Value = 0.003128754952689148

# Row 5
$ws.Range("D5").Value = 0.0001894159094018535
$ws.Range("E5").Value = 0.0001894159094018535

# Row 6
$ws.Range("D6").Value = 0.7391647162181709
$ws.Range("E6").Value = 0.7391647162181709

# Row 7
$ws.Range("D7").Value = 0.9999999683112414
$ws.Range("E7").Value = [double]"3.168875861714326E-08"

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.2251247457340769
$ws.Range("E8").Value = 0.7748752542659232

# Row 9
$ws.Range("D9").Value = 0.5791340029517763
$ws.Range("E9").Value = 0.4208659970482237

# Row 10
$ws.Range("D10").Value = 0.8823237379389625
$ws.Range("E10").Value = 0.1176762620610375

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.0008103353063904546
$ws.Range("E11").Value = 0.9991896646936096
$ws.Range("F11").Value = 1.351793169975281
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.9980750372295302
$ws.Range("E12").Value = 0.9980750372295302

# Row 13
$ws.Range("D13").Value = [double]"2.389245021323969E-06"
$ws.Range("E13").Value = [double]"2.389245021323969E-06"

# Row 14
$ws.Range("D14").Value = 0.001406661838687808
$ws.Range("E14").Value = 0.001406661838687808

# Row 15
$ws.Range("D15").Value = [double]"3.97337176711579E-05"
$ws.Range("E15").Value = [double]"3.97337176711579E-05"

# Row 16
$ws.Range("D16").Value = 0.7708014978983186
$ws.Range("E16").Value = 0.7708014978983186

# Row 17
$ws.Range("D17").Value = 0.9999999995207094
$ws.Range("E17").Value = [double]"4.792906072026426E-10"

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = 0.004079516649088826
$ws.Range("E18").Value = 0.9959204833509112

# Row 19
$ws.Range("D19").Value = 0.9769829433387397
$ws.Range("E19").Value = 0.0230170566612603

# Row 20
$ws.Range("D20").Value = 0.6341044172542758
$ws.Range("E20").Value = 0.3658955827457242

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = [double]"1.084996148672812E-08"
$ws.Range("E21").Value = 0.9999999891500385
$ws.Range("F21").Value = 3.20471715927124
$ws.Range("G21").Value = 0.6
